$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Climbing Stairs" row (row 18): entered Name, Link, Description, Approach
# (this order controls the shared-string table insertion order)
$ws.Range("A18").Value = "Climbing Stairs"
$ws.Range("D18").Value = "https://leetcode.com/problems/climbing-stairs/"
$ws.Range("B18").Value = "Return dinstics way to climb for given steps"
$ws.Range("C18").Value = "Use a bottom-up approach using for loop and an array. Subproblem: dp[n] = dp[n-1] + dp[n-2]"

# Add the external hyperlink on D18 (same pattern as the other rows)
$ws.Hyperlinks.Add($ws.Range("D18"), "https://leetcode.com/problems/climbing-stairs/")

# Copy cell formatting from the row above so the new row matches the sheet's
# existing look (A -> "Good" style, B/C -> normal, D -> "Hyperlink" style)
# without introducing brand-new style entries.
$ws.Range("A14").Copy()
$ws.Range("A18").PasteSpecial(-4122)
$ws.Range("B14").Copy()
$ws.Range("B18").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("C18").PasteSpecial(-4122)
$ws.Range("D14").Copy()
$ws.Range("D18").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Move the active selection like in the authored workbook
$ws.Range("C19").Select() | Out-Null
